$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "None" placeholder text from B43 (Dumper ID unknown) -> empty cell
$ws.Range("B43").Value = ""

# Add the new row 44 data: rear-left-outer tyre reading
$ws.Range("A44").Value = "2024-12-13 04:33:38"
$ws.Range("B44").Value = "None"
$ws.Range("C44").Value = "rear-left-outer"
$ws.Range("D44").Value = 1920
$ws.Range("E44").Value = 1920
